# Update values in Sheet1 to reflect new RandomForest imputation results.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("E7").Value = 16.23510000000001
$ws.Range("C9").Value = -10.441
$ws.Range("E12").Value = 18.10590000000002
$ws.Range("C18").Value = -12.44619999999999
$ws.Range("C20").Value = -11.55350000000001
$ws.Range("E26").Value = 16.17029999999999
$ws.Range("C27").Value = -12.28349999999999
$ws.Range("E27").Value = 16.63809999999999
$ws.Range("E29").Value = 17.06330000000002
$ws.Range("E37").Value = 16.63200000000002
$ws.Range("E38").Value = 16.3798
$ws.Range("E51").Value = 16.98640000000001
$ws.Range("E55").Value = 16.3151
$ws.Range("C69").Value = -11.2931
$ws.Range("E69").Value = 17.12120000000003
$ws.Range("E70").Value = 17.93740000000002
$ws.Range("C76").Value = -12.3144
$ws.Range("C82").Value = -11.93529999999999
$ws.Range("E83").Value = 16.6586
$ws.Range("E102").Value = 16.8146
